# Adds the "withdrawal use case" admin-activity rows (36-50) to the "2022_8" log sheet,
# matching the upstream data export growth from A1:N35 to A1:N50.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2022_8")

# Row 36
$ws.Cells.Item(36, 1).Value = "Tue Aug 23 2022"
$ws.Cells.Item(36, 2).Value = "16:08:34 GMT+0000 (Coordinated Universal Time)"
$ws.Cells.Item(36, 4).Value = "User"
$ws.Cells.Item(36, 5).Value = "/api/auth/login"
$ws.Cells.Item(36, 6).Value = "login"
$ws.Cells.Item(36, 7).Value = "succeeded"
$ws.Cells.Item(36, 8).Value = "developer@nex-softwares.com  login"
$ws.Cells.Item(36, 14).Value = "developer@nex-softwares.com"

# Row 37
$ws.Cells.Item(37, 1).Value = "Tue Aug 23 2022"
$ws.Cells.Item(37, 2).Value = "16:08:35 GMT+0000 (Coordinated Universal Time)"
$ws.Cells.Item(37, 4).Value = "User"
$ws.Cells.Item(37, 5).Value = "/api/user"
$ws.Cells.Item(37, 6).Value = "read"
$ws.Cells.Item(37, 7).Value = "succeeded"
$ws.Cells.Item(37, 8).Value = "NEX  Admin  read all users (4) from 0 to 100"
$ws.Cells.Item(37, 10).Value = 9
$ws.Cells.Item(37, 11).Value = "all"
$ws.Cells.Item(37, 12).Value = "NEX"
$ws.Cells.Item(37, 13).Value = "Admin"

# Row 38
$ws.Cells.Item(38, 1).Value = "Tue Aug 23 2022"
$ws.Cells.Item(38, 2).Value = "16:08:38 GMT+0000 (Coordinated Universal Time)"
$ws.Cells.Item(38, 4).Value = "User"
$ws.Cells.Item(38, 5).Value = "/api/user"
$ws.Cells.Item(38, 6).Value = "read"
$ws.Cells.Item(38, 7).Value = "succeeded"
$ws.Cells.Item(38, 8).Value = "NEX  Admin  read all users (4) from 0 to 100"
$ws.Cells.Item(38, 10).Value = 9
$ws.Cells.Item(38, 11).Value = "all"
$ws.Cells.Item(38, 12).Value = "NEX"
$ws.Cells.Item(38, 13).Value = "Admin"

# Row 39
$ws.Cells.Item(39, 1).Value = "Tue Aug 23 2022"
$ws.Cells.Item(39, 2).Value = "16:08:53 GMT+0000 (Coordinated Universal Time)"
$ws.Cells.Item(39, 4).Value = "User"
$ws.Cells.Item(39, 5).Value = "/api/to-validate-users"
$ws.Cells.Item(39, 6).Value = "read"
$ws.Cells.Item(39, 7).Value = "succeeded"
$ws.Cells.Item(39, 8).Value = "NEX  Admin  read all to validate users (1) from 0 to 100"
$ws.Cells.Item(39, 10).Value = 9
$ws.Cells.Item(39, 11).Value = "all"
$ws.Cells.Item(39, 12).Value = "NEX"
$ws.Cells.Item(39, 13).Value = "Admin"

# Row 40
$ws.Cells.Item(40, 1).Value = "Tue Aug 23 2022"
$ws.Cells.Item(40, 2).Value = "16:09:11 GMT+0000 (Coordinated Universal Time)"
$ws.Cells.Item(40, 4).Value = "Trip"
$ws.Cells.Item(40, 5).Value = "/api/trip"
$ws.Cells.Item(40, 6).Value = "read"
$ws.Cells.Item(40, 7).Value = "succeeded"
$ws.Cells.Item(40, 8).Value = "NEX  Admin  read all trips (2) from 0 to 100"
$ws.Cells.Item(40, 10).Value = 9
$ws.Cells.Item(40, 11).Value = "all"
$ws.Cells.Item(40, 12).Value = "NEX"
$ws.Cells.Item(40, 13).Value = "Admin"

# Row 41
$ws.Cells.Item(41, 1).Value = "Tue Aug 23 2022"
$ws.Cells.Item(41, 2).Value = "16:09:18 GMT+0000 (Coordinated Universal Time)"
$ws.Cells.Item(41, 4).Value = "Preference"
$ws.Cells.Item(41, 5).Value = "/api/preference"
$ws.Cells.Item(41, 6).Value = "read"
$ws.Cells.Item(41, 7).Value = "succeeded"
$ws.Cells.Item(41, 8).Value = "NEX  Admin  read all preferences (undefined) from undefined to NaN"
$ws.Cells.Item(41, 10).Value = 9
$ws.Cells.Item(41, 11).Value = "all"
$ws.Cells.Item(41, 12).Value = "NEX"
$ws.Cells.Item(41, 13).Value = "Admin"

# Row 42
$ws.Cells.Item(42, 1).Value = "Tue Aug 23 2022"
$ws.Cells.Item(42, 2).Value = "16:09:22 GMT+0000 (Coordinated Universal Time)"
$ws.Cells.Item(42, 4).Value = "VehicleType"
$ws.Cells.Item(42, 5).Value = "/api/vehicle-type"
$ws.Cells.Item(42, 6).Value = "read"
$ws.Cells.Item(42, 7).Value = "succeeded"
$ws.Cells.Item(42, 8).Value = "NEX  Admin  read all vehicle types (undefined) from undefined to NaN"
$ws.Cells.Item(42, 10).Value = 9
$ws.Cells.Item(42, 11).Value = "all"
$ws.Cells.Item(42, 12).Value = "NEX"
$ws.Cells.Item(42, 13).Value = "Admin"

# Row 43
$ws.Cells.Item(43, 1).Value = "Tue Aug 23 2022"
$ws.Cells.Item(43, 2).Value = "16:09:27 GMT+0000 (Coordinated Universal Time)"
$ws.Cells.Item(43, 4).Value = "Pricing"
$ws.Cells.Item(43, 5).Value = "/api/pricing"
$ws.Cells.Item(43, 6).Value = "read"
$ws.Cells.Item(43, 7).Value = "succeeded"
$ws.Cells.Item(43, 8).Value = "NEX  Admin  read all pricing (undefined) from undefined to NaN"
$ws.Cells.Item(43, 10).Value = 9
$ws.Cells.Item(43, 11).Value = "all"
$ws.Cells.Item(43, 12).Value = "NEX"
$ws.Cells.Item(43, 13).Value = "Admin"

# Row 44
$ws.Cells.Item(44, 1).Value = "Tue Aug 23 2022"
$ws.Cells.Item(44, 2).Value = "16:09:51 GMT+0000 (Coordinated Universal Time)"
$ws.Cells.Item(44, 4).Value = "VehicleType"
$ws.Cells.Item(44, 5).Value = "/api/vehicle-type"
$ws.Cells.Item(44, 6).Value = "read"
$ws.Cells.Item(44, 7).Value = "succeeded"
$ws.Cells.Item(44, 8).Value = "NEX  Admin  read all vehicle types (undefined) from undefined to NaN"
$ws.Cells.Item(44, 10).Value = 9
$ws.Cells.Item(44, 11).Value = "all"
$ws.Cells.Item(44, 12).Value = "NEX"
$ws.Cells.Item(44, 13).Value = "Admin"

# Row 45
$ws.Cells.Item(45, 1).Value = "Tue Aug 23 2022"
$ws.Cells.Item(45, 2).Value = "16:09:55 GMT+0000 (Coordinated Universal Time)"
$ws.Cells.Item(45, 4).Value = "Trip"
$ws.Cells.Item(45, 5).Value = "/api/trip"
$ws.Cells.Item(45, 6).Value = "read"
$ws.Cells.Item(45, 7).Value = "succeeded"
$ws.Cells.Item(45, 8).Value = "NEX  Admin  read all trips (2) from 0 to 100"
$ws.Cells.Item(45, 10).Value = 9
$ws.Cells.Item(45, 11).Value = "all"
$ws.Cells.Item(45, 12).Value = "NEX"
$ws.Cells.Item(45, 13).Value = "Admin"

# Row 46
$ws.Cells.Item(46, 1).Value = "Tue Aug 23 2022"
$ws.Cells.Item(46, 2).Value = "16:41:59 GMT+0000 (Coordinated Universal Time)"
$ws.Cells.Item(46, 4).Value = "Trip"
$ws.Cells.Item(46, 5).Value = "/api/trip"
$ws.Cells.Item(46, 6).Value = "read"
$ws.Cells.Item(46, 7).Value = "succeeded"
$ws.Cells.Item(46, 8).Value = "NEX  Admin  read all trips (2) from 0 to 100"
$ws.Cells.Item(46, 10).Value = 9
$ws.Cells.Item(46, 11).Value = "all"
$ws.Cells.Item(46, 12).Value = "NEX"
$ws.Cells.Item(46, 13).Value = "Admin"

# Row 47
$ws.Cells.Item(47, 1).Value = "Tue Aug 23 2022"
$ws.Cells.Item(47, 2).Value = "21:29:18 GMT+0000 (Coordinated Universal Time)"
$ws.Cells.Item(47, 4).Value = "User"
$ws.Cells.Item(47, 5).Value = "/api/user"
$ws.Cells.Item(47, 6).Value = "read"
$ws.Cells.Item(47, 7).Value = "succeeded"
$ws.Cells.Item(47, 8).Value = "NEX  Admin  read all users (4) from 0 to 100"
$ws.Cells.Item(47, 10).Value = 9
$ws.Cells.Item(47, 11).Value = "all"
$ws.Cells.Item(47, 12).Value = "NEX"
$ws.Cells.Item(47, 13).Value = "Admin"

# Row 48
$ws.Cells.Item(48, 1).Value = "Tue Aug 23 2022"
$ws.Cells.Item(48, 2).Value = "21:38:30 GMT+0000 (Coordinated Universal Time)"
$ws.Cells.Item(48, 4).Value = "Trip"
$ws.Cells.Item(48, 5).Value = "/api/trip"
$ws.Cells.Item(48, 6).Value = "read"
$ws.Cells.Item(48, 7).Value = "succeeded"
$ws.Cells.Item(48, 8).Value = "NEX  Admin  read all trips (2) from 0 to 100"
$ws.Cells.Item(48, 10).Value = 9
$ws.Cells.Item(48, 11).Value = "all"
$ws.Cells.Item(48, 12).Value = "NEX"
$ws.Cells.Item(48, 13).Value = "Admin"

# Row 49
$ws.Cells.Item(49, 1).Value = "Tue Aug 23 2022"
$ws.Cells.Item(49, 2).Value = "21:38:39 GMT+0000 (Coordinated Universal Time)"
$ws.Cells.Item(49, 4).Value = "Trip"
$ws.Cells.Item(49, 5).Value = "/api/trip"
$ws.Cells.Item(49, 6).Value = "read"
$ws.Cells.Item(49, 7).Value = "succeeded"
$ws.Cells.Item(49, 8).Value = "NEX  Admin  read all trips (2) from 0 to 100"
$ws.Cells.Item(49, 10).Value = 9
$ws.Cells.Item(49, 11).Value = "all"
$ws.Cells.Item(49, 12).Value = "NEX"
$ws.Cells.Item(49, 13).Value = "Admin"

# Row 50
$ws.Cells.Item(50, 1).Value = "Tue Aug 23 2022"
$ws.Cells.Item(50, 2).Value = "21:38:40 GMT+0000 (Coordinated Universal Time)"
$ws.Cells.Item(50, 4).Value = "User"
$ws.Cells.Item(50, 5).Value = "/api/to-validate-users"
$ws.Cells.Item(50, 6).Value = "read"
$ws.Cells.Item(50, 7).Value = "succeeded"
$ws.Cells.Item(50, 8).Value = "NEX  Admin  read all to validate users (1) from 0 to 100"
$ws.Cells.Item(50, 10).Value = 9
$ws.Cells.Item(50, 11).Value = "all"
$ws.Cells.Item(50, 12).Value = "NEX"
$ws.Cells.Item(50, 13).Value = "Admin"
